$wb = $excel.ActiveWorkbook

# --- Sheet "pro" (sheet1): update B2:B26 values ---
$wsPro = $wb.Worksheets.Item("pro")
$proValues = @{
    2 = 337599.23039294971
    3 = 362951.61224775022
    4 = 420763.40510366077
    5 = 465567.04254975641
    6 = 453489.69997539261
    7 = 460677.98622973036
    8 = 461019.05086613767
    9 = 451628.46852432966
    10 = 443452.94885510672
    11 = 445530.54834986478
    12 = 473256.60256800201
    13 = 488502.65745268983
    14 = 566501.38070431782
    15 = 557991.36183137295
    16 = 567088.49908260431
    17 = 512928.31637944595
    18 = 561927.28138269507
    19 = 624842.41767691402
    20 = 670833.67550894385
    21 = 780135
    22 = 898015
    23 = 995874.9786596772
    24 = 1006122.1998036421
    25 = 1048336.1663152121
    26 = 1025395.0151284173
}
foreach ($row in $proValues.Keys) {
    $wsPro.Range("B$row").Value = $proValues[$row]
}

# --- Sheet "ind" (sheet2): update B2:B101 values ---
$wsInd = $wb.Worksheets.Item("ind")
$indValues = @{
    2 = 413211.22011070466
    3 = 370293.89622573956
    4 = 351990.82376866246
    5 = 391149.4417952309
    6 = 410608.95664496807
    7 = 392913.98064049159
    8 = 388809.68591225473
    9 = 448957.85618627351
    10 = 471293.37803646922
    11 = 455317.90897936106
    12 = 449884.52878473717
    13 = 526223.30037304363
    14 = 532918.04767048755
    15 = 506558.28820240102
    16 = 510417.86876110855
    17 = 555429.83489742293
    18 = 506745.29963468225
    19 = 506267.30208682985
    20 = 505841.23790225614
    21 = 531855.68184337753
    22 = 492920.01441526582
    23 = 506027.54263660102
    24 = 490676.47100788995
    25 = 593591.38508607924
    26 = 493332.934580828
    27 = 506478.95245021739
    28 = 518119.30627533549
    29 = 566826.53600971028
    30 = 513432.27079214522
    31 = 506596.49094508
    32 = 495484.87585829542
    33 = 526779.27582395379
    34 = 497096.65052396245
    35 = 485441.68499434448
    36 = 477981.20719741873
    37 = 544803.14617708966
    38 = 481062.31012012897
    39 = 497343.50856755173
    40 = 499218.22058495588
    41 = 537093.68784470041
    42 = 506856.73993156833
    43 = 513126.75387143268
    44 = 508043.02192936128
    45 = 612070.20760238287
    46 = 489048.41517508513
    47 = 509119.92162051989
    48 = 557417.65664146538
    49 = 653454.36834974203
    50 = 603515.40325159091
    51 = 625376.92122822965
    52 = 637155.73239790439
    53 = 695707.53624306025
    54 = 633290.68025018054
    55 = 610449.66476757277
    56 = 621609.66835897358
    57 = 657922.72835501947
    58 = 600703.36325582722
    59 = 610295.12165164249
    60 = 645334.72252770641
    61 = 708077.37266468722
    62 = 598431.70010082272
    63 = 586764.8678945438
    64 = 551846.33342811931
    65 = 582451.85730316024
    66 = 589326.40187105082
    67 = 606429.89028324082
    68 = 619190.59444716328
    69 = 656396.54812031344
    70 = 648863.67060696986
    71 = 668568.04468012624
    72 = 688299.2292029435
    73 = 758400.21436383796
    74 = 755186.03399985679
    75 = 806334.72610893392
    76 = 810088.30736297264
    77 = 884308.59401328303
    78 = 917904.3291761179
    79 = 963848.58719728806
    80 = 996713.74927677389
    81 = 1121533.3343498204
    82 = 1143221.9553929029
    83 = 1155155.1343545543
    84 = 1175350.2343760885
    85 = 1290955.0599179533
    86 = 1328485.0204159489
    87 = 1358206.138994155
    88 = 1365045.5390068183
    89 = 1442132.2300571376
    90 = 1487433.3372979718
    91 = 1492827.7652067312
    92 = 1502356.260133937
    93 = 1599221.3944183504
    94 = 1673718.0049533916
    95 = 1654225.4405021081
    96 = 1632235.3734587387
    97 = 1714193.9940752485
    98 = 1709005.7404231201
    99 = 1738274.1056583151
    100 = 1746707.6879806812
    101 = 1813857.5273328633
}
foreach ($row in $indValues.Keys) {
    $wsInd.Range("B$row").Value = $indValues[$row]
}

# --- Sheet "conso" (sheet4): update B2:B26 values ---
$wsConso = $wb.Worksheets.Item("conso")
$consoValues = @{
    2 = 117932.70064230703
    3 = 136717.73945925164
    4 = 151276.07595356726
    5 = 147350.80719245897
    6 = 149686.25039106773
    7 = 149796.89633306218
    8 = 146745.71935722887
    9 = 144089.67305308857
    10 = 144764.94602217639
    11 = 153773.49517346997
    12 = 158727.93463469512
    13 = 184071.97420141255
    14 = 181306.84572758398
    15 = 184262.729920365
    16 = 166664.38007658831
    17 = 182585.23852771285
    18 = 203028.23648071659
    19 = 217971.44817132014
    20 = 253486
    21 = 335766
    22 = 323998.02481358481
    23 = 464205.39584834577
    24 = 483682.10656863073
    25 = 473097.50146801601
    26 = 3489936.780383375
}
foreach ($row in $consoValues.Keys) {
    $wsConso.Range("B$row").Value = $consoValues[$row]
}

# --- Recalculate workbook so VA sheet (formulas) pick up new cached values ---
$excel.CalculateFull()

# --- Adjust column B width on "pro" sheet (15.6328125 -> ~14.6328125) ---
$wsPro.Columns.Item(2).ColumnWidth = 13.83

# --- Update selections/view to match target ---
$wsPro.Activate()
$wsPro.Range("B2:B101").Select()

$wsInd.Activate()
$wsInd.Range("B2:B101").Select()

$wsVA = $wb.Worksheets.Item("VA")
$wsVA.Activate()
$wsVA.Range("B2:B101").Select()

$wsConso.Activate()
$wsConso.Range("B2:B101").Select()

$wsPro.Activate()
